$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header in A1 was renamed from "Country" to "NAME"
$ws.Range("A1").Value = "NAME"
